# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
#
# - Two pairs of columns (B and F:AC) get corrected on five match rows
#   where the wrong fixture's data had been attached to the row
#   (the two games in each pair had their results/odds swapped).
# - The two most recent (not yet played / placeholder) fixtures at the
#   bottom of the sheet are removed.

function ToRow($vals) {
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the 5 rows whose match data was swapped with its neighbour ---

$ws.Range("B5").Value = 6221786
$ws.Range("F5:AC5").Value = (ToRow @("Ordabasy", "FK Atyrau", 2, 1, "H", 1.571, 3.8, 4.75, 1.5, 4, 5.25, -1, 1.875, 1.925, 2.5, 1.9, 1.9, 0.5, -1, -1, 0, 0, 0.8999999999999999, -1))

$ws.Range("B6").Value = 6726054
$ws.Range("F6:AC6").Value = (ToRow @("FK Aktobe", "Shakhter Karagandy", 2, 0, "H", 1.727, 4, 3.5, 1.333, 4.75, 6.5, -1.5, 1.975, 1.825, 3, 1.975, 1.825, 0.333, -1, -1, 0.9750000000000001, -1, -1, 0.825))

$ws.Range("B16").Value = 6221693
$ws.Range("F16:AC16").Value = (ToRow @("Zhetysu", "Shakhter Karagandy", 1, 3, "A", 2, 3.4, 3.1, 2.2, 3.3, 2.8, -0.25, 1.95, 1.85, 2.5, 1.85, 1.95, -1, -1, 1.8, -1, 0.8500000000000001, 0.8500000000000001, -1))

$ws.Range("B17").Value = 6221698
$ws.Range("F17:AC17").Value = (ToRow @("FK Maktaaral", "FK Aktobe", 1, 2, "A", 4.333, 3.5, 1.666, 4.2, 3.4, 1.7, 0.75, 1.825, 1.975, 2.5, 1.925, 1.875, -1, -1, 0.7, -0.5, 0.4875, 0.925, -1))

$ws.Range("B37").Value = 6221712
$ws.Range("F37:AC37").Value = (ToRow @("FK Aksu", "Shakhter Karagandy", 2, 1, "H", 2.1, 3.25, 3, 2.15, 3.25, 2.9, -0.25, 1.95, 1.85, 2.5, 1.975, 1.825, 1.15, -1, -1, 0.95, -1, 0.9750000000000001, -1))

$ws.Range("B38").Value = 6221708
$ws.Range("F38:AC38").Value = (ToRow @("Kaisar Kyzylorda", "Kairat Almaty", 0, 0, "D", 3, 3.4, 2.05, 3.2, 3.4, 1.95, 0.5, 1.75, 1.95, 2.25, 1.925, 1.875, -1, 2.4, -1, 0.75, -1, -1, 0.875))

$ws.Range("B99").Value = 6221753
$ws.Range("F99:AC99").Value = (ToRow @("FK Aksu", "Tobol Kostanay", 0, 3, "A", 2.75, 3.1, 2.375, 2.625, 3.2, 2.45, 0, 2, 1.8, 2.5, 1.9, 1.9, -1, -1, 1.45, -1, 0.8, 0.8999999999999999, -1))

$ws.Range("B100").Value = 6221752
$ws.Range("F100:AC100").Value = (ToRow @("FK Kyzylzhar", "Kaisar Kyzylorda", 0, 1, "A", 1.833, 3.2, 4, 1.85, 3.2, 4, -0.5, 1.9, 1.9, 2, 1.775, 2.025, -1, -1, 3, -1, 0.8999999999999999, -1, 1.025))

$ws.Range("B119").Value = 7873759
$ws.Range("F119:AC119").Value = (ToRow @("FK Aktobe", "FK Zhenys", 3, 0, "H", 1.25, 5.75, 7, 1.444, 4.75, 4.75, -1.25, 1.95, 1.85, 2.75, 1.9, 1.9, 0.444, -1, -1, 0.95, -1, 0.45, -0.5))

$ws.Range("B120").Value = 7874795
$ws.Range("F120:AC120").Value = (ToRow @("FK Kyzylzhar", "Tobol Kostanay", 0, 0, "D", 2.2, 3.1, 3, 2.625, 3, 2.55, 0, 1.9, 1.9, 2, 1.95, 1.85, -1, 2, -1, 0, 0, -1, 0.8500000000000001))

# --- Remove the two trailing placeholder fixtures (rows 129 & 130) ---

$ws.Range("A129:A130").EntireRow.Delete()
